$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.79
$ws.Range("I2").Value = 4.7

# Row 3
$ws.Range("G3").Value = 1.47
$ws.Range("H3").Value = 8.2
$ws.Range("I3").Value = 8.6
$ws.Range("N3").Value = 5.4
$ws.Range("T3").Value = 1.89
$ws.Range("U3").Value = 2.06
$ws.Range("V3").Value = 1.13
$ws.Range("W3").Value = 3.1
$ws.Range("AN3").Value = 5.7

# Row 4
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 3.75
$ws.Range("L4").Value = 1.32
$ws.Range("N4").Value = 5.1
$ws.Range("Q4").Value = 1.66
$ws.Range("R4").Value = 1.57
$ws.Range("V4").Value = 1.36
$ws.Range("AO4").Value = 28

# Row 5
$ws.Range("F5").Value = 1.52
$ws.Range("G5").Value = 1.72
$ws.Range("H5").Value = 4.5
$ws.Range("I5").Value = 6.4
$ws.Range("J5").Value = 4.3
$ws.Range("K5").Value = 5.5
$ws.Range("L5").Value = 1.22
$ws.Range("O5").Value = 1.16
$ws.Range("P5").Value = 2.6
$ws.Range("Q5").Value = 1.45
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 2.12
$ws.Range("T5").Value = 1.56
$ws.Range("U5").Value = 2.3
$ws.Range("V5").Value = 1.18
$ws.Range("W5").Value = 2.38
$ws.Range("X5").Value = 32
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 55
$ws.Range("AB5").Value = 15.5
$ws.Range("AC5").Value = 13.5
$ws.Range("AD5").Value = 25
$ws.Range("AE5").Value = 75
$ws.Range("AF5").Value = 14.5
$ws.Range("AG5").Value = 11.5
$ws.Range("AH5").Value = 19
$ws.Range("AI5").Value = 60
$ws.Range("AJ5").Value = 19
$ws.Range("AK5").Value = 16
$ws.Range("AL5").Value = 27
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 6.6
$ws.Range("AO5").Value = 1000

# Row 6
$ws.Range("G6").Value = 3
$ws.Range("I6").Value = 2.64
$ws.Range("K6").Value = 4.2
$ws.Range("L6").Value = 1.29
$ws.Range("U6").Value = 2.52
$ws.Range("V6").Value = 1.61

# Row 7
$ws.Range("G7").Value = 1.6
$ws.Range("I7").Value = 9.2
$ws.Range("N7").Value = 4.3
$ws.Range("S7").Value = 2.6
$ws.Range("U7").Value = 1.96
$ws.Range("W7").Value = 2.66
$ws.Range("AH7").Value = 27

# Row 8
$ws.Range("F8").Value = 2.48
$ws.Range("G8").Value = 2.52
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 3.35
$ws.Range("W8").Value = 1.66
$ws.Range("AN8").Value = 25

# Row 9
$ws.Range("F9").Value = 2.6
$ws.Range("G9").Value = 2.64
$ws.Range("H9").Value = 2.68
$ws.Range("I9").Value = 2.7
$ws.Range("Q9").Value = 1.53
$ws.Range("R9").Value = 1.75
$ws.Range("S9").Value = 2.26
$ws.Range("U9").Value = 2.92
$ws.Range("V9").Value = 1.58

# Row 10
$ws.Range("F10").Value = 9.8
$ws.Range("K10").Value = 5.8
$ws.Range("Q10").Value = 1.64
$ws.Range("T10").Value = 1.97
$ws.Range("AO10").Value = 5.1

# Row 11
$ws.Range("P11").Value = 2.94
$ws.Range("AF11").Value = 8
$ws.Range("AN11").Value = 3.35

# Row 12
$ws.Range("I12").Value = 11.5
$ws.Range("W12").Value = 4.3

# Row 13
$ws.Range("F13").Value = 6
$ws.Range("G13").Value = 6.2
$ws.Range("H13").Value = 1.65
$ws.Range("I13").Value = 1.67
$ws.Range("V13").Value = 2.5
$ws.Range("W13").Value = 1.19
$ws.Range("Y13").Value = 9.2
$ws.Range("AE13").Value = 16.5

# Row 14
$ws.Range("H14").Value = 2.4
$ws.Range("I14").Value = 2.42
$ws.Range("K14").Value = 3.8
$ws.Range("N14").Value = 4.8
$ws.Range("Q14").Value = 1.75
$ws.Range("T14").Value = 1.62
$ws.Range("V14").Value = 1.7
$ws.Range("AH14").Value = 15

# Row 15
$ws.Range("F15").Value = 2.42
$ws.Range("G15").Value = 3.1
$ws.Range("H15").Value = 2.34
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 3.45
$ws.Range("L15").Value = 1.26
$ws.Range("M15").Value = 1.03
$ws.Range("Q15").Value = 1.53
$ws.Range("R15").Value = 1.49
$ws.Range("T15").Value = 1.5
$ws.Range("V15").Value = 1.5
$ws.Range("W15").Value = 1.48
$ws.Range("X15").Value = 26
$ws.Range("Y15").Value = 17.5
$ws.Range("AA15").Value = 1000
$ws.Range("AB15").Value = 17
$ws.Range("AC15").Value = 10.5
$ws.Range("AF15").Value = 24
$ws.Range("AG15").Value = 14.5
$ws.Range("AH15").Value = 17
$ws.Range("AJ15").Value = 1000
$ws.Range("AK15").Value = 30
$ws.Range("AL15").Value = 36
$ws.Range("AM15").Value = 70
$ws.Range("AN15").Value = 19.5
$ws.Range("AO15").Value = 18.5

# Row 16
$ws.Range("F16").Value = 2.52
$ws.Range("G16").Value = 2.88
$ws.Range("I16").Value = 3.2
$ws.Range("L16").Value = 1.42
$ws.Range("T16").Value = 1.75
$ws.Range("V16").Value = 1.46
$ws.Range("Y16").Value = 13.5
$ws.Range("Z16").Value = 22
$ws.Range("AC16").Value = 9
$ws.Range("AF16").Value = 20
$ws.Range("AH16").Value = 18.5
$ws.Range("AO16").Value = 1000
